$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.136.65'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.892.71'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7390'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.80'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9996'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3175'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.46%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07214'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.96'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08342'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.88%  '
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7610'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.459'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.77%  '
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.829.77'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.18'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.172'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.076.11'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '250.56'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.66'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007891'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.188.30'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9981'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.958'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1583'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.319'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.74'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.78'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.066'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.485'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.590'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.53%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05369'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.06%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7796'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.21%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.732'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01966'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.765'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4577'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.098.39'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.081'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.60'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8726'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '104.74'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.13%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.868'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.605'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.636'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.060.05'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.47%  '
